# correção da tabela de monitorias
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the now-unused placeholder rows (4:10) - only 2 data rows remain
$ws.Rows("4:10").Delete() | Out-Null

# 2) Clear the old "left aligned" placeholder style from the remaining data rows (2:3)
#    so they fall back to the default (unstyled) cell format
$ws.Range("B2:J3").ClearFormats() | Out-Null

# 3) Re-enter the header row, reordering the columns
$ws.Range("A1").Value = "data"
$ws.Range("B1").Value = "projeto"
$ws.Range("C1").Value = "nome_analista"
$ws.Range("D1").Value = "id_atendimento"
$ws.Range("E1").Value = "duracao"
$ws.Range("F1").Value = "chamado"
$ws.Range("G1").Value = "nome_cliente"
$ws.Range("H1").Value = "categoria"
$ws.Range("I1").Value = "nota"
$ws.Range("J1").Value = "observacao"

# 4) Row 2 data
#    A2 looks like a date to Excel's auto-detection, so force it to stay plain text
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2024-04-23"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "Flowserve"
$ws.Range("C2").Value = "Teste"
$ws.Range("D2").Value = 123
$ws.Range("E2").Value = 321
$ws.Range("F2").Value = 123234
$ws.Range("G2").Value = "Teste"
$ws.Range("H2").Value = "Acessos"
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = "abc"

# 5) Row 3 data
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2024-04-23"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = "Cteep"

#    C3:G3, I3 and J3 hold the text "234" (not the number 234), so force text too
$ws.Range("C3:G3").NumberFormat = "@"
$ws.Range("C3").Value = "234"
$ws.Range("D3").Value = "234"
$ws.Range("E3").Value = "234"
$ws.Range("F3").Value = "234"
$ws.Range("G3").Value = "234"
$ws.Range("C3:G3").Style = "Normal"

$ws.Range("H3").Value = "Acessos"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "234"
$ws.Range("I3").Style = "Normal"

$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "234"
$ws.Range("J3").Style = "Normal"

# 6) Restore the active selection seen in the saved workbook
$ws.Range("F8").Select() | Out-Null
